# Falta ver vista previa encuestas
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Criterios" (sheet3): add helper column C (LEFT(A,3)), apply wrap-text
# style to column B, highlight rows 5-6 (A5:B6) in yellow, widen/narrow
# column B, adjust the sheet's zoom + selection + print orientation.
# ---------------------------------------------------------------------------
$wsCriterios = $wb.Worksheets.Item("Criterios")

# New helper column: idCriterio truncated to the "x.y" characteristic code
$wsCriterios.Range("C2:C54").Formula = "=LEFT(A2,3)"

# Highlight the two rows that still need review (A5:B6) in yellow
$wsCriterios.Range("A5:B6").Interior.Color = 65535

# Column B: narrower, word-wrapped
$wsCriterios.Columns.Item(2).ColumnWidth = 40
$wsCriterios.Range("B1:B54").WrapText = $true

# Row heights follow straight from Excel's word-wrap autofit once the
# narrower, wrapped column B is in place.
$wsCriterios.Rows.Item("1:54").AutoFit()

# Print the sheet in portrait orientation
$wsCriterios.PageSetup.Orientation = 1

# Zoom in on the sheet and select A2:A6 (first block of criteria)
$wsCriterios.Application.ActiveWindow.Zoom = 160
$wsCriterios.Range("A2:A6").Select()

# ---------------------------------------------------------------------------
# Sheet "Caracteristicas" (sheet2): move the working selection here (the
# characteristics still pending review).
# ---------------------------------------------------------------------------
$wsCaracteristicas = $wb.Worksheets.Item("Caracteristicas")
$r1 = $wsCaracteristicas.Range("A10")
$r2 = $wsCaracteristicas.Range("A11")
$r3 = $wsCaracteristicas.Range("A12")
$r4 = $wsCaracteristicas.Range("A14")
$r5 = $wsCaracteristicas.Range("A13")
$union = $excel.Union($r1, $r2, $r3, $r4, $r5)
$union.Select()
